# Applies the cryptos list price/volume update described in the commit
# "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.893.01"
$ws.Range("E2").Value = "  -0.75%  "
$ws.Range("D3").Value = "2.037.49"
$ws.Range("E3").Value = "  -1.11%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.38"
$ws.Range("E5").Value = "  -0.87%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.611"
$ws.Range("E6").Value = "  -0.64%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.56"
$ws.Range("E7").Value = "  +1.40%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.380"
$ws.Range("E9").Value = "  -1.88%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0821"
$ws.Range("E10").Value = "  +1.23%  "
$ws.Range("E11").Value = "  -0.15%  "
$ws.Range("D12").Value = "2.340.26"
$ws.Range("E12").Value = "  -1.08%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.53"
$ws.Range("E13").Value = "  -1.68%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.22"
$ws.Range("E14").Value = "  -0.56%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.762"
$ws.Range("E15").Value = "  +0.75%  "
$ws.Range("E16").Value = "  -2.35%  "
$ws.Range("D17").Value = "2.035.17"
$ws.Range("E17").Value = "  -1.26%  "
$ws.Range("D18").Value = "37.844.19"
$ws.Range("E18").Value = "  -0.64%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "69.88"
$ws.Range("E19").Value = "  -0.10%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.91"
$ws.Range("E20").Value = "  -5.90%  "
$ws.Range("D21").Value = "0.0₃0825"
$ws.Range("E21").Value = "  -1.41%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "224.02"
$ws.Range("E22").Value = "  -0.61%  "
$ws.Range("E23").Value = "  +0.10%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.42"
$ws.Range("E24").Value = "  -0.39%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.27"
$ws.Range("E25").Value = "  +0.68%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.38"
$ws.Range("E26").Value = "  +0.53%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "167.34"
$ws.Range("E27").Value = "  +0.69%  "
$ws.Range("E28").Value = "  -2.17%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.88"
$ws.Range("E29").Value = "  -0.92%  "
$ws.Range("E30").Value = "  -3.92%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.121"
$ws.Range("E31").Value = "  +0.92%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.22"
$ws.Range("E32").Value = "  +7.88%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.42"
$ws.Range("E33").Value = "  -3.38%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0607"
$ws.Range("E34").Value = "  -0.21%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.52"
$ws.Range("E35").Value = "  -1.91%  "
$ws.Range("E36").Value = "  +3.70%  "
$ws.Range("E37").Value = "  -2.76%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.33"
$ws.Range("E38").Value = "  +1.14%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "17.73"
$ws.Range("E40").Value = "  +4.31%  "
$ws.Range("D41").Value = "1.540.96"
$ws.Range("E41").Value = "  +0.36%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0218"
$ws.Range("E42").Value = "  +0.08%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "96.56"
$ws.Range("E43").Value = "  -1.95%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.82"
$ws.Range("E44").Value = "  -2.00%  "
$ws.Range("E45").Value = "  -1.12%  "
$ws.Range("E46").Value = "  -2.91%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.01"
$ws.Range("E47").Value = "  -2.55%  "
$ws.Range("E48").Value = "  -1.18%  "
$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.19"
$ws.Range("E49").Value = "  +0.91%  "
$ws.Range("B50").Value = "MXToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.96"
$ws.Range("E50").Value = "  -0.49%  "
$ws.Range("D51").Value = "2.228.98"
$ws.Range("E51").Value = "  -1.05%  "
